$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.947.80'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '1.637.68'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.82%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.80'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.255'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.95%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0637'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.61'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').Value = '1.864.34'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('D14').Value = '1.608.07'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.52'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '25.953.57'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '193.65'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.38'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.93'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('E23').Value = '  -1.62%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '144.23'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('E27').Value = '  +2.49%  '
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.49'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0501'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('E32').Value = '  -1.22%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('E34').Value = '  -2.88%  '
$ws.Range('E35').Value = '  +1.72%  '
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').Value = '1.139.02'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('E40').Value = '  +0.61%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '99.27'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.799'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.41'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.67%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.773.65'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0114'
$ws.Range('E45').Value = '  +10.37%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '56.46'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0529'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.46'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.66'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.415'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.02%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0962'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.19%  '
